$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the ages between Ram (row 2) and the row-3 person,
# and rename "Shyam" to "Hari" in row 3.
$ws.Range("C2").Value = 20
$ws.Range("B3").Value = "Hari"
$ws.Range("C3").Value = 23
